# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# and the priority ("ht" -> "mt") to reflect a freshly generated handback report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-31 13:09:19"
$wsOverview.Range("G3").Value = "2016-08-31 13:09:19"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-31 13:09:02"
$wsZhCn.Range("H3").Value = "2016-08-31 13:09:02"
$wsZhCn.Range("K2").Value = "2016-08-31 13:09:57"
$wsZhCn.Range("K3").Value = "2016-08-31 13:09:57"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-31 13:09:19"
$wsDeDe.Range("H3").Value = "2016-08-31 13:09:19"
$wsDeDe.Range("K2").Value = "2016-08-31 13:10:35"
$wsDeDe.Range("K3").Value = "2016-08-31 13:10:35"
